$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.109.76'
$ws.Range('E2').Value = '  -0.23%  '

# Row 3
$ws.Range('D3').Value = '3.238.05'
$ws.Range('E3').Value = '  +0.13%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '529.12'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.38%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '171.19'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.54%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.596'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.34%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').Value = '3.237.71'
$ws.Range('E9').Value = '  +0.30%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.605'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.14%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.06'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -5.87%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.133'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.42%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000254'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.87%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.12'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.35%  '

# Row 15
$ws.Range('D15').Value = '3.763.60'
$ws.Range('E15').Value = '  +0.25%  '

# Row 16
$ws.Range('E16').Value = '  -0.79%  '

# Row 17
$ws.Range('D17').Value = '3.242.15'
$ws.Range('E17').Value = '  +0.39%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '63.085.18'
$ws.Range('E18').Value = '  +0.10%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '17.21'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.24%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.06'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.69%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.967'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.64%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '366.32'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.91%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.75'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.40%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '80.96'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.55%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.17'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.34%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.97'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +8.19%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.08'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.39%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.64'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.90%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.24'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.33%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.18'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.20%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '28.44'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.29%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '633.53'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.41%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.45'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.79%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.19'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.85%  '

# Row 35
$ws.Range('E35').Value = '  +4.91%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '56.58'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.22%  '

# Row 37
$ws.Range('E37').Value = '  -0.08%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '36.58'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.50%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.377'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.44%  '

# Row 40
$ws.Range('D40').Value = '0.0₃0716'
$ws.Range('E40').Value = '  +12.58%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.17%  '

# Row 42
$ws.Range('E42').Value = '  +2.05%  '

# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.882.86'
$ws.Range('E43').Value = '  +1.98%  '

# Row 44
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.54'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +10.88%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.94'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +6.84%  '

# Row 46
$ws.Range('E46').Value = '  +4.92%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0393'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +5.21%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.09'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.61%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.58'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.86%  '

# Row 50
$ws.Range('E50').Value = '  +3.16%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '134.13'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.50%  '
